$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Day 4 rows (48-57): fill in "Original File Name" (column A) and the
#     newly-counted slip stats (columns F and G); row 48 is also flagged
#     Questionable (column H) ---
# F56 alone uses the alternate "s=5" style (matching cells like A39), all
# other F/G cells use the regular "s=2" style (matching cells like C48).

$day4 = @(
    @{ Row = 48; A = "MVI_0629.MOV"; F = 21.0; FStyle = 2; G = 117.0; Questionable = $true  },
    @{ Row = 49; A = "MVI_0630.MOV"; F = 15.0; FStyle = 2; G = 52.0;  Questionable = $false },
    @{ Row = 50; A = "MVI_0631.MOV"; F = 5.0;  FStyle = 2; G = 30.0;  Questionable = $false },
    @{ Row = 51; A = "MVI_0632.MOV"; F = 3.0;  FStyle = 2; G = 26.0;  Questionable = $false },
    @{ Row = 52; A = "MVI_0633.MOV"; F = 2.0;  FStyle = 2; G = 21.0;  Questionable = $false },
    @{ Row = 53; A = "MVI_0634.MOV"; F = 4.0;  FStyle = 2; G = 16.0;  Questionable = $false },
    @{ Row = 54; A = "MVI_0635.MOV"; F = 1.0;  FStyle = 2; G = 15.0;  Questionable = $false },
    @{ Row = 55; A = "MVI_0636.MOV"; F = 1.0;  FStyle = 2; G = 13.0;  Questionable = $false },
    @{ Row = 56; A = "MVI_0637.MOV"; F = 1.5;  FStyle = 5; G = 12.0;  Questionable = $false },
    @{ Row = 57; A = "MVI_0638.MOV"; F = 2.0;  FStyle = 2; G = 11.0;  Questionable = $false }
)

# Reference cells already carrying the styles we need to reproduce:
#   style "5" (used by column A here, and by F56) -> e.g. A39
#   style "2" (used by column B/F/G/H here)       -> e.g. C48
$styleRefS5 = $ws.Cells.Item(39, 1)
$styleRefS2 = $ws.Cells.Item(48, 3)

foreach ($entry in $day4) {
    $r = $entry.Row

    # Column A: Original File Name
    $ws.Cells.Item($r, 1).Value = $entry.A
    $styleRefS5.Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)

    # Column F: Crop begin time
    $ws.Cells.Item($r, 6).Value = $entry.F
    if ($entry.FStyle -eq 5) {
        $styleRefS5.Copy()
    } else {
        $styleRefS2.Copy()
    }
    $ws.Cells.Item($r, 6).PasteSpecial(-4122)

    # Column G: Crop end time
    $ws.Cells.Item($r, 7).Value = $entry.G
    $styleRefS2.Copy()
    $ws.Cells.Item($r, 7).PasteSpecial(-4122)

    # Column H: Questionable
    if ($entry.Questionable) {
        $ws.Cells.Item($r, 8).Value = "Yes"
        $styleRefS2.Copy()
        $ws.Cells.Item($r, 8).PasteSpecial(-4122)
    }
}

# --- Existing rows also newly flagged Questionable ("Yes" in column H) ---
$questionableRows = @(38, 58, 59)
foreach ($r in $questionableRows) {
    $ws.Cells.Item($r, 8).Value = "Yes"
    $styleRefS2.Copy()
    $ws.Cells.Item($r, 8).PasteSpecial(-4122)
}
